# Fruta / hortaliza, semanal
# Insert a new week of price data (4 quality rows: Especial, Primera, Segunda, Tercera)
# at the top of the data block (row 22), pushing the existing history down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 22; existing rows 22:45 become 26:49.
$ws.Rows("22:25").Insert()

# Common (constant) metadata values shared by every data row in this sheet.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$productoId = 100101
$producto  = "Berries"
$categoriaId = 100112025
$categoria = "Frutilla"
$variedad  = "Sin especificar"
$unidad    = "$/bandeja 3 kilos"
$origen    = "Región de Arica y Parinacota"
$kgUnidad  = 3

# New week: 2022-11-29 (serial date 44894)
$fecha = 44894

# Row 22: Especial
$r = 22
$ws.Range("A$r").Value = $mercadoId
$ws.Range("B$r").Value = $mercado
$ws.Range("C$r").Value = $region
$ws.Range("D$r").Value = $fecha
$ws.Range("E$r").Value = $codreg
$ws.Range("F$r").Value = $tipo
$ws.Range("G$r").Value = $productoId
$ws.Range("H$r").Value = $producto
$ws.Range("I$r").Value = $categoriaId
$ws.Range("J$r").Value = $categoria
$ws.Range("K$r").Value = $variedad
$ws.Range("L$r").Value = "Especial"
$ws.Range("M$r").Value = 60
$ws.Range("N$r").Value = 7000
$ws.Range("O$r").Value = 8000
$ws.Range("P$r").Value = 7500
$ws.Range("Q$r").Value = $unidad
$ws.Range("R$r").Value = $origen
$ws.Range("S$r").Value = 2500
$ws.Range("T$r").Value = $kgUnidad

# Row 23: Primera
$r = 23
$ws.Range("A$r").Value = $mercadoId
$ws.Range("B$r").Value = $mercado
$ws.Range("C$r").Value = $region
$ws.Range("D$r").Value = $fecha
$ws.Range("E$r").Value = $codreg
$ws.Range("F$r").Value = $tipo
$ws.Range("G$r").Value = $productoId
$ws.Range("H$r").Value = $producto
$ws.Range("I$r").Value = $categoriaId
$ws.Range("J$r").Value = $categoria
$ws.Range("K$r").Value = $variedad
$ws.Range("L$r").Value = "Primera"
$ws.Range("M$r").Value = 70
$ws.Range("N$r").Value = 6000
$ws.Range("O$r").Value = 7000
$ws.Range("P$r").Value = 6500
$ws.Range("Q$r").Value = $unidad
$ws.Range("R$r").Value = $origen
$ws.Range("S$r").Value = 2167
$ws.Range("T$r").Value = $kgUnidad

# Row 24: Segunda
$r = 24
$ws.Range("A$r").Value = $mercadoId
$ws.Range("B$r").Value = $mercado
$ws.Range("C$r").Value = $region
$ws.Range("D$r").Value = $fecha
$ws.Range("E$r").Value = $codreg
$ws.Range("F$r").Value = $tipo
$ws.Range("G$r").Value = $productoId
$ws.Range("H$r").Value = $producto
$ws.Range("I$r").Value = $categoriaId
$ws.Range("J$r").Value = $categoria
$ws.Range("K$r").Value = $variedad
$ws.Range("L$r").Value = "Segunda"
$ws.Range("M$r").Value = 72
$ws.Range("N$r").Value = 5000
$ws.Range("O$r").Value = 6000
$ws.Range("P$r").Value = 5486
$ws.Range("Q$r").Value = $unidad
$ws.Range("R$r").Value = $origen
$ws.Range("S$r").Value = 1829
$ws.Range("T$r").Value = $kgUnidad

# Row 25: Tercera
$r = 25
$ws.Range("A$r").Value = $mercadoId
$ws.Range("B$r").Value = $mercado
$ws.Range("C$r").Value = $region
$ws.Range("D$r").Value = $fecha
$ws.Range("E$r").Value = $codreg
$ws.Range("F$r").Value = $tipo
$ws.Range("G$r").Value = $productoId
$ws.Range("H$r").Value = $producto
$ws.Range("I$r").Value = $categoriaId
$ws.Range("J$r").Value = $categoria
$ws.Range("K$r").Value = $variedad
$ws.Range("L$r").Value = "Tercera"
$ws.Range("M$r").Value = 74
$ws.Range("N$r").Value = 4000
$ws.Range("O$r").Value = 5000
$ws.Range("P$r").Value = 4500
$ws.Range("Q$r").Value = $unidad
$ws.Range("R$r").Value = $origen
$ws.Range("S$r").Value = 1500
$ws.Range("T$r").Value = $kgUnidad

# Make sure date cells use the same date/time number format as the rest of column D.
$ws.Range("D22:D25").NumberFormat = $ws.Range("D26").NumberFormat
